# Applies the "Fixed AI_Summary, refactor project structure and cleanup" edit:
#  - Add a new "Email" row (A7/B7) with a mailto: hyperlink on B7
#  - Rename the "NrRooms" label (A5) to "Rooms"
#  - Update MinPrice (B2) and MaxPrice (B3) values
#  - Move the active selection to D3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new Email row first, including the hyperlink, so the new shared
# strings ("Email" / the address) are appended before we touch A5.
$ws.Range("A7").Value = "Email"
$ws.Range("B7").Value = "tunaru.alexandra2005@gmail.com"
$ws.Hyperlinks.Add($ws.Range("B7"), "mailto:tunaru.alexandra2005@gmail.com") | Out-Null

# Rename NrRooms -> Rooms
$ws.Range("A5").Value = "Rooms"

# Update MinPrice / MaxPrice values
$ws.Range("B2").Value = 70000
$ws.Range("B3").Value = 130000

# Match the new active cell selection
$ws.Range("D3").Select() | Out-Null
